$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"-0.05479815075587835"
$ws.Range("C2").Value = [double]"-0.043967052216493"
$ws.Range("D2").Value = [double]"0.8023455465196794"
$ws.Range("B3").Value = [double]"-0.043967052216493"
$ws.Range("C3").Value = [double]"0.0005085818044777701"
$ws.Range("D3").Value = [double]"0.01156733915145191"
$ws.Range("B4").Value = [double]"0.0005085818044777701"
$ws.Range("C4").Value = [double]"-3.758094060735129e-05"
$ws.Range("D4").Value = [double]"0.07389360035391107"
$ws.Range("B5").Value = [double]"-3.758094060735129e-05"
$ws.Range("C5").Value = [double]"-2.61509441301655e-06"
$ws.Range("D5").Value = [double]"0.06958565620640708"
$ws.Range("B6").Value = [double]"-2.61509441301655e-06"
$ws.Range("C6").Value = [double]"2.662474662962211e-07"
$ws.Range("D6").Value = [double]"0.101811798828747"
$ws.Range("B7").Value = [double]"2.662474662962211e-07"
$ws.Range("C7").Value = [double]"5.168764061203035e-08"
$ws.Range("D7").Value = [double]"0.1941338309470476"
$ws.Range("B8").Value = [double]"5.168764061203035e-08"
$ws.Range("C8").Value = [double]"-3.258101455294593e-09"
$ws.Range("D8").Value = [double]"0.06303443950460116"
$ws.Range("B9").Value = [double]"-3.258101455294593e-09"
$ws.Range("C9").Value = [double]"-1.521550663241555e-09"
$ws.Range("D9").Value = [double]"0.4670053048130077"
$ws.Range("B10").Value = [double]"-1.521550663241555e-09"
$ws.Range("C10").Value = [double]"1.757871626040242e-11"
$ws.Range("D10").Value = [double]"0.01155315868546382"
$ws.Range("B11").Value = [double]"1.757871626040242e-11"
$ws.Range("C11").Value = [double]"5.363054444984527e-11"
$ws.Range("D11").Value = [double]"3.05087946442669"
$ws.Range("B12").Value = [double]"5.363054444984527e-11"
$ws.Range("C12").Value = [double]"2.524758180300068e-12"
$ws.Range("D12").Value = [double]"0.04707687020893842"
$ws.Range("B13").Value = [double]"2.524758180300068e-12"
$ws.Range("C13").Value = [double]"-2.049249658853114e-12"
$ws.Range("D13").Value = [double]"0.8116617562991952"
$ws.Range("B14").Value = [double]"-2.049249658853114e-12"
$ws.Range("C14").Value = [double]"-2.410294186461215e-13"
$ws.Range("D14").Value = [double]"0.117618376855564"
$ws.Range("B15").Value = [double]"-2.410294186461215e-13"
$ws.Range("C15").Value = [double]"7.949196856316121e-14"
$ws.Range("D15").Value = [double]"0.3298019345923537"
$ws.Range("B16").Value = [double]"7.949196856316121e-14"
$ws.Range("C16").Value = [double]"1.720845688168993e-14"
$ws.Range("D16").Value = [double]"0.2164804469273743"
$ws.Range("B17").Value = [double]"1.720845688168993e-14"
$ws.Range("C17").Value = [double]"-2.220446049250313e-15"
$ws.Range("D17").Value = [double]"0.1290322580645161"
$ws.Range("B18").Value = [double]"-2.220446049250313e-15"
$ws.Range("C18").Value = [double]"-1.998401444325282e-15"
$ws.Range("D18").Value = [double]"0.9"
$ws.Range("B19").Value = [double]"-1.998401444325282e-15"
$ws.Range("C19").Value = [double]"1.332267629550188e-15"
$ws.Range("D19").Value = [double]"0.6666666666666666"
